# Apply updated "想去人数" (want-to-go count) values to the
# "展览" (Exhibition) and "全部类型" (All Types) sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F9").Value  = 6899
$wsExhibition.Range("F16").Value = 16366
$wsExhibition.Range("F17").Value = 5
$wsExhibition.Range("F25").Value = 1093
$wsExhibition.Range("F26").Value = 4510
$wsExhibition.Range("F27").Value = 367

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 6899
$wsAll.Range("F18").Value = 16366
$wsAll.Range("F19").Value = 5
$wsAll.Range("F29").Value = 1093
$wsAll.Range("F30").Value = 4510
$wsAll.Range("F31").Value = 367
